$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9872071743011475
$ws.Range("B1").Value = 1.984674572944641
$ws.Range("C1").Value = 8.639707565307617
$ws.Range("D1").Value = 2.797600746154785
$ws.Range("E1").Value = 1.424680471420288
